$wb = $excel.ActiveWorkbook

# --- Data sheet: add F3 / F4 values ---
$wsData = $wb.Worksheets.Item("Data")
$wsData.Range("F3").Value = 123
$wsData.Range("F4").Value = 1443

# --- Define sheet: B5 text changes from "object" to "class" ---
$wsDefine = $wb.Worksheets.Item("Define")
$wsDefine.Range("B5").Value = "class"

# --- Update selections on both sheets ---
# Define sheet is no longer the active/selected tab; just move its selection.
$wsDefine.Range("F6").Select()

# Data sheet becomes the active tab with F8 selected.
$wsData.Activate()
$wsData.Range("F8").Select()
